# Daily update at 8 AM UTC
# Appends the next day's results to the "Wins Over Time" tracking sheet,
# and flips the previous last row's date format back to the regular
# (non-final) style while the newly appended row becomes the "latest" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row that used to be "last" (row 5) loses its special last-row
# date formatting and reverts to the standard date/time format used by
# every other non-final row (A2:A4).
$ws.Range("A5").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's row of data.
$ws.Range("A6").Value = 45955
$ws.Range("A6").NumberFormat = "YYYY-MM-DD"
$ws.Range("B6").Value = 9
$ws.Range("C6").Value = 14
$ws.Range("D6").Value = 10
